$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dimDates")

# --- Row 7: this used to be the "latest" plain-link entry (2024-04-18).
# It is now superseded, so it becomes an archived entry: a real date value
# (styled like the other archived rows) plus wayback-machine URLs.
$ws.Range("A7").Value = 45412
$ws.Range("A7").NumberFormat = $ws.Range("A2").NumberFormat
$ws.Range("B7").Value = "https://web.archive.org/web/20240430111825/https://www.oryxspioenkop.com/2022/02/attack-on-europe-documenting-equipment.html"
$ws.Range("C7").Value = "https://web.archive.org/web/20240430110305/https://www.oryxspioenkop.com/2022/02/attack-on-europe-documenting-ukrainian.html"

# --- Row 8 (new): the new "latest" entry for 2024-05-01, with plain
# (non-archived) links, matching the pattern the old row 7 used to follow.
# Force the date cell to be stored as literal text (not auto-converted to a
# date serial) by temporarily marking it as Text before assigning the value,
# then restore the default style so the cell carries no explicit format.
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "2024-05-01"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").Value = "https://www.oryxspioenkop.com/2022/02/attack-on-europe-documenting-equipment.html"
$ws.Range("C8").Value = "https://www.oryxspioenkop.com/2022/02/attack-on-europe-documenting-ukrainian.html"

$wb.Save()
